# Append the new wishlist entry as the next row in the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "La dependienta"
$ws.Range("B21").Value = "Sayaka Murata"

# Touch C21 (left blank, same as the other rows in this column) so the
# row's used range/dimension picks it up without pulling in a new style.
$ws.Range("C21").Style = "Normal"
